$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cell F1 with the same style as the other headers (B1:E1)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate F2:F72 with the recorded time_taken timestamps
$timestamps = @(
    "2021-10-05 10:51:13.571520",
    "2021-10-05 10:51:13.571532",
    "2021-10-05 10:51:13.571536",
    "2021-10-05 10:51:13.571538",
    "2021-10-05 10:51:13.571541",
    "2021-10-05 10:51:13.571544",
    "2021-10-05 10:51:13.571546",
    "2021-10-05 10:51:13.571549",
    "2021-10-05 10:51:13.571551",
    "2021-10-05 10:51:13.571554",
    "2021-10-05 10:51:13.571557",
    "2021-10-05 10:51:13.571559",
    "2021-10-05 10:51:13.571561",
    "2021-10-05 10:51:13.571564",
    "2021-10-05 10:51:13.571566",
    "2021-10-05 10:51:13.571569",
    "2021-10-05 10:51:13.571572",
    "2021-10-05 10:51:13.571574",
    "2021-10-05 10:51:13.571577",
    "2021-10-05 10:51:13.571579",
    "2021-10-05 10:51:13.571582",
    "2021-10-05 10:51:13.571585",
    "2021-10-05 10:51:13.571587",
    "2021-10-05 10:51:13.571590",
    "2021-10-05 10:51:13.571592",
    "2021-10-05 10:51:13.571595",
    "2021-10-05 10:51:13.571597",
    "2021-10-05 10:51:13.571600",
    "2021-10-05 10:51:13.571602",
    "2021-10-05 10:51:13.571605",
    "2021-10-05 10:51:13.571607",
    "2021-10-05 10:51:13.571610",
    "2021-10-05 10:51:13.571612",
    "2021-10-05 10:51:13.571615",
    "2021-10-05 10:51:13.571617",
    "2021-10-05 10:51:13.571620",
    "2021-10-05 10:51:13.571622",
    "2021-10-05 10:51:13.571625",
    "2021-10-05 10:51:13.571627",
    "2021-10-05 10:51:13.571630",
    "2021-10-05 10:51:13.571633",
    "2021-10-05 10:51:13.571635",
    "2021-10-05 10:51:13.571638",
    "2021-10-05 10:51:13.571640",
    "2021-10-05 10:51:13.571643",
    "2021-10-05 10:51:13.571645",
    "2021-10-05 10:51:13.571648",
    "2021-10-05 10:51:13.571650",
    "2021-10-05 10:51:13.571653",
    "2021-10-05 10:51:13.571655",
    "2021-10-05 10:51:13.571658",
    "2021-10-05 10:51:13.571660",
    "2021-10-05 10:51:13.571663",
    "2021-10-05 10:51:13.571665",
    "2021-10-05 10:51:13.571668",
    "2021-10-05 10:51:13.571670",
    "2021-10-05 10:51:13.571673",
    "2021-10-05 10:51:13.571675",
    "2021-10-05 10:51:13.571678",
    "2021-10-05 10:51:13.571680",
    "2021-10-05 10:51:13.571682",
    "2021-10-05 10:51:13.571685",
    "2021-10-05 10:51:13.571687",
    "2021-10-05 10:51:13.571690",
    "2021-10-05 10:51:13.571693",
    "2021-10-05 10:51:13.571696",
    "2021-10-05 10:51:13.571699",
    "2021-10-05 10:51:13.571701",
    "2021-10-05 10:51:13.571704",
    "2021-10-05 10:51:13.571706",
    "2021-10-05 10:51:13.571709"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
